# Weekly refresh of the "Hortaliza, Feria Lagunitas de Puerto Montt - Apio"
# price sheet: a new weekly record is inserted as row 263 (pushing the
# existing rows 263:276 down to 264:277), growing the sheet from 276 to
# 277 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 263:276 down to 264:277, leaving a blank row 263 to fill in.
$ws.Rows.Item(263).EntireRow.Insert()

# Populate the new row 263 with the latest weekly price entry.
$ws.Cells.Item(263, 1).Value  = 4
$ws.Cells.Item(263, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(263, 3).Value  = "Los Lagos"
$ws.Cells.Item(263, 4).Value  = 44753
$ws.Cells.Item(263, 5).Value  = 10
$ws.Cells.Item(263, 6).Value  = 100112017
$ws.Cells.Item(263, 7).Value  = "Apio"
$ws.Cells.Item(263, 8).Value  = "Americana (o)"
$ws.Cells.Item(263, 9).Value  = "Primera"
$ws.Cells.Item(263, 10).Value = 25
$ws.Cells.Item(263, 11).Value = 11000
$ws.Cells.Item(263, 12).Value = 11000
$ws.Cells.Item(263, 13).Value = 11000
$ws.Cells.Item(263, 14).Value = "$/docena de matas"
$ws.Cells.Item(263, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(263, 16).Value = 1833
$ws.Cells.Item(263, 17).Value = 6
$ws.Cells.Item(263, 18).Value = "Hortaliza"
